$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume figures per latest data refresh
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.793.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.637.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "322.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.89"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.048.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.647.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.784.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.89"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.94"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.48%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.10"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.50"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.47"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.30"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0797"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.94"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.72%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.17"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.11"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.63"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.111"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0315"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.072.69"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.18"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.68"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.16"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.70%  "
